$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DANH SÁCH NỢ")

# Remove the existing "Điều khoản dịch vụ" hyperlink before the row shift so it
# isn't left dangling on the old cell reference.
$ws.Range("D16").Hyperlinks.Delete()

# The row for "Trần Huỳnh Như Ý / Mua thẻ Mobi" (already paid in full / "Đã trả đủ")
# was removed from the debt list; Excel shifts every row below it up by one.
$ws.Rows("7:7").Delete()

# Re-create the hyperlink on its new location (one row up) and restore the
# large hyperlink font that the cell had before (Add() resets it to the
# default 11pt hyperlink font).
$ws.Hyperlinks.Add($ws.Range("D15"), "https://tinyurl.com/dieukhoan29")
$ws.Range("D15").Font.Size = 18

# Reapply the AutoFilter over the now-smaller data range.
$ws.AutoFilterMode = $false
$ws.Range("A1:M18").AutoFilter()

# Keep the workbook-level _xlnm._FilterDatabase name in sync with the filter range.
$name = $wb.Names.Item(1)
$name.RefersTo = "='DANH SÁCH NỢ'!`$A`$1:`$M`$18"

# Restore the last-used selection.
$ws.Range("L21").Select()
